$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Rename the worksheet tab ----
$ws.Name = "ExploradorEntidades"

# ---- Row 1: drop the fixed 60pt height back to the sheet default ----
$ws.Rows(1).AutoFit()

# ---- Row 2: taller custom height ----
$ws.Rows(2).RowHeight = 87.75

# ---- Grab the "blank, no-wrap, vertically centered" look (currently only ----
# ---- on C4) for the two cells that will need it later, before C4 itself ----
# ---- is restyled.                                                       ----
$ws.Range("C4").Copy()
$ws.Range("E3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C4").Copy()
$ws.Range("F4").PasteSpecial(-4122)

# ---- D2: "Positivo" -> "eCenter" (format/style untouched) ----
$ws.Range("D2").Value = "eCenter"

# ---- Row 3 reshaping: copy matching column formats from row 2 / row 3/4 ----
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("H3").PasteSpecial(-4122)

$ws.Range("B3").Value = "Selección de elemento secundario(ONT)" + [char]10
$ws.Range("D3").Value = "eCenter"
$ws.Range("E3").Value = "Debe haber accedido a la vista"
$ws.Range("F3").Value = "1.Clic en ""elemento secundario""" + [char]10 + "2.Clic en tarjeta ONT"
$ws.Range("G3").Value = "N/A"
$ws.Range("H3").Value = "El sistema debe redirigido correctamente a los elementos secundarios (ONT)"

# ---- Row 4 gets a new custom height + reshaped content ----
$ws.Rows(4).RowHeight = 66.75

$ws.Range("A4").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("A4").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("E4").PasteSpecial(-4122)

$ws.Range("A4").Value = "CP_EXPENT_003"
$ws.Range("B4").Value = "Crear nuevo registro" + [char]10 + "entidad(ONT)"
$ws.Range("C4").Value = "Positivo"
$ws.Range("D4").Value = "eCenter"
$ws.Range("E4").Value = "Debe haber accedido al apartado de elemento secundario(ONT)"
$ws.Range("F4").ClearContents()

# ---- Active selection moves from C4 to F4 ----
$null = $ws.Range("F4").Select()
